# Translate English facilitator-guide text to Italian.
#
# We locate each piece of English text with Find.Execute (no replacement,
# just a match) and then assign the Italian text directly to the found
# Range's .Text property. Doing the substitution this way (rather than
# passing ReplaceWith to Find.Execute) avoids Word's AutoCorrect/AutoFormat
# "smart quotes" feature silently turning straight apostrophes into curly
# ones, which would not match the target OOXML.
#
# wdFindContinue = 1 (used for Forward scanning, no wrap)

$d = $word.ActiveDocument

$script:cursor = 0

function Replace-Next($find, $replace) {
    $rng = $d.Content
    $rng.Start = $script:cursor
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $find"
    }
    $rng.Text = $replace
    $script:cursor = $rng.End
}

$pairs = @(
    @("Video Title", "Titolo del Video"),
    @("The Airport Problem", "Il Problema dell'Aeroporto"),
    @("Topic", "Argomento"),
    @("Geometry", "Geometria"),
    @("Aim(s)", "Obiettivo/i"),
    @("Get the intuitive idea of a minimization problem, figure out how to practically implement minimization problems.", "Ottenere l'idea intuitiva di un problema di minimizzazione, capire come implementare praticamente i problemi di minimizzazione."),
    @("Length", "Lunghezza"),
    @("Camp Location", "Posizione del Campo"),
    @("Facilitators", "Mediatori"),
    @("N. of students", "N. di studenti"),
    @("Date", "Data"),
    @("Resources", "Risorse"),
    @("needed", "necessarie"),
    @("Pins (3 each group), string (1/group), metal ring (optional but convenient to avoid friction 1/group), thick cardboard or wooden disposable surface (1/group)", "Spille (3 per gruppo), filo (1/gruppo), anello metallico (facoltativo ma conveniente per evitare frizione 1/gruppo), cartone spesso o superficie usa e getta di legno (1/gruppo)"),
    @("Preparations", "Preparazioni"),
    @("Pin 3 points on the wood", "Spilla 3 punti sul legno"),
    @("Video time", "Tempo del video"),
    @("What facilitator does", "Cosa fa il facilitatore"),
    @("What learners do", "Cosa fanno gli studenti"),
    @("General VMC Video Introduction", "Introduzione Generale al Video di VMC"),
    @("Video Introduction", "Introduzione al video"),
    @("Riddle", "Enigma"),
    @("Introduction of the first experiment", "Introduzione del primo esperimento"),
    @("VIDEO PAUSE", "PAUSA VIDEO"),
    @("Finding a solution", "Trovare una soluzione"),
    @("Assist the process, provoke thoughts", "Assiste il processo, causa pensieri"),
    @("Try to find a setting of the string such that the minimization of the length of the string corresponds to minimize the sum of the lengths of the roads", "Prova a trovare un'impostazione del filo così che la minimizzazione della lunghezza del filo corrisponda a minimizzare la somma delle lunghezze delle strade"),
    @("Solution ", "Soluzione "),
    @("VIDEO PAUSE", "PAUSA VIDEO"),
    @("Geometry", "Geometria"),
    @("Assist the process, provoke thoughts", "Assite il processo, causa pensieri"),
    @("Try to figure out what geometrical property the new point has in relation to the starting 3.", "Prova a capire che proprietà geometrica ha il nuovo punto in relazione alla partenza 3."),
    @("Showing the 120° angles", "Mostra gli angoli a 120°"),
    @("solution", "soluzione"),
    @("The point ,F, found as shown in the video, is called Fermat Point.", "Il punto F, trovato come mostrato nel video è detto Punto di Fermat."),
    @("There are several different possible paths of the string that can be used to find point F.", "Esistono molti percorsi diversi del filo, utilizzabili per trovare il punto F."),
    @("Notice that the use of the ring is not strictly necessary, but it helps to reduce friction (natural enemy of this experience).", "Nota che l'uso dell'anello non è strettamente necessario, ma aiuta a ridurre la frizione (nemico naturale di quest'esperienza)."),
    @("Once the point is found (before watching the solution in the video) students can be asked to find the angles by noticing that each angle is congruent and they form 360° all together.", "Una volta trovato il punto (prima di guardare la soluzione nel video), gli studenti possono esser invitati a trovare gli angoli notando che ognuno è congruente e che formano, tutti insieme, 360°."),
    @("A geometrical construction that can be used to find F consists in building equilateral triangles on the sides of the original triangle and connecting opposite points:", "Una costruzione geometrica utilizzabile per trovare F consiste nel costruire triangoli equilateri sui lati del triangolo originale e connettendo i punti opposti:"),
    @("This construction can be replicated on the wooden board to verify that the two methods will lead to the same point.", "Questa costruzione è replicabile sulla scheda di legno per verificare che i due metodi condurranno allo stesso punto.")
)

foreach ($pair in $pairs) {
    Replace-Next $pair[0] $pair[1]
}

Write-Output "Done: applied $($pairs.Count) replacements"
